# "Se marcan en rojo las areas por completar"
#
# 1) Slide 3 ("Descripcion de ZOFRI S.A."): the << ... >> placeholder
#    bullets get red text.
# 2) Slide 5 ("Estructura Organizacional de ZOFRI S.A."): the
#    Directores/Gerentes/Dotacion box is moved, a new "Ejecutivos" line
#    is added between Gerentes and Dotacion, and every line is turned red.

$p = $ppt.ActivePresentation

$red = 255  # RGB(255,0,0) -> pure red, stored as <a:srgbClr val="FF0000"/>

# --- Slide 3: "<< Resena breve >>" / "<< Principal actividad >>" / "<< Ubicacion >>" ---
$s3 = $p.Slides.Item(3)
$ph3 = $s3.Shapes.Placeholders.Item(2)
$tr3 = $ph3.TextFrame.TextRange
$count3 = $tr3.Paragraphs().Count
for ($i = 1; $i -le $count3; $i++) {
    $tr3.Paragraphs($i, 1).Font.Color.RGB = $red
}

# --- Slide 5: Directores / Gerentes / (new) Ejecutivos / Dotacion ---
$s5 = $p.Slides.Item(5)
$ph5 = $s5.Shapes.Placeholders.Item(2)

# Reposition the box.
$ph5.Left = 605.6543307086614
$ph5.Top = 179.99787401574804

$tr5 = $ph5.TextFrame.TextRange

# "Dotacion" is currently the 3rd paragraph (Directores, Gerentes, Dotacion).
# Insert the new "Ejecutivos" paragraph right before it.
$dotacionPara = $tr5.Paragraphs(3, 1)
$dotacionPara.InsertBefore("Ejecutivos" + [char]13)

# Colour every line (Directores, Gerentes, Ejecutivos, Dotacion) red.
$count5 = $tr5.Paragraphs().Count
for ($i = 1; $i -le $count5; $i++) {
    $tr5.Paragraphs($i, 1).Font.Color.RGB = $red
}
